$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows 3-13 (weekly "Haba" price records) were re-shuffled: each
# row's record-specific fields (Fecha, Volumen, Precio mínimo/máximo/
# promedio, Origen, Precio $/Kg) now carry the data that used to live in
# a different row. Apply the full target values for each row.

# Row 3
$ws.Range("D3").Value = 44435
$ws.Range("J3").Value = 15
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = 14000
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 560

# Row 4
$ws.Range("D4").Value = 44435
$ws.Range("J4").Value = 15
$ws.Range("K4").Value = 14000
$ws.Range("L4").Value = 14000
$ws.Range("M4").Value = 14000
$ws.Range("O4").Value = "Provincia del Elquí"
$ws.Range("P4").Value = 560

# Row 5
$ws.Range("D5").Value = 44340
$ws.Range("J5").Value = 25
$ws.Range("K5").Value = 15000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 15000
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 600

# Row 6
$ws.Range("D6").Value = 44421
$ws.Range("J6").Value = 20
$ws.Range("K6").Value = 15000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 15000
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 600

# Row 7
$ws.Range("D7").Value = 44446
$ws.Range("J7").Value = 15
$ws.Range("K7").Value = 13000
$ws.Range("L7").Value = 13000
$ws.Range("M7").Value = 13000
$ws.Range("O7").Value = "Provincia de Limarí"
$ws.Range("P7").Value = 520

# Row 8
$ws.Range("D8").Value = 44453
$ws.Range("J8").Value = 55
$ws.Range("K8").Value = 14000
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 14455
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 578

# Row 9
$ws.Range("D9").Value = 44425
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 14000
$ws.Range("L9").Value = 14000
$ws.Range("M9").Value = 14000
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 560

# Row 10
$ws.Range("D10").Value = 44376
$ws.Range("J10").Value = 15
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 12000
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 480

# Row 11
$ws.Range("D11").Value = 44418
$ws.Range("J11").Value = 12
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = 15000
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 600

# Row 12
$ws.Range("D12").Value = 44432
$ws.Range("J12").Value = 15
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 14000
$ws.Range("M12").Value = 14000
$ws.Range("O12").Value = "Provincia del Elquí"
$ws.Range("P12").Value = 560

# Row 13
$ws.Range("D13").Value = 44467
$ws.Range("J13").Value = 35
$ws.Range("K13").Value = 12000
$ws.Range("L13").Value = 12000
$ws.Range("M13").Value = 12000
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 480
